$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The module list in E6 used to point at "Blueoth HC05"; replace it with the
# new "StepMotor" entry (driving a stepper motor directly, without a driver
# board) per "dua len mach dieu khien dong co buoc khong co Driver".
$ws.Range("E6").Value = "StepMotor"

# Move the current selection to D6.
$ws.Range("D6").Select()
